$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.362.66"
$ws.Range("E2").Value = "  +3.99%  "

$ws.Range("D3").Value = "1.806.35"
$ws.Range("E3").Value = "  +1.87%  "

$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "316.04"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "0.5509"
$ws.Range("E7").Value = "  +5.69%  "

$ws.Range("D8").Value = "0.3856"
$ws.Range("E8").Value = "  +6.62%  "

$ws.Range("D9").Value = "0.07592"
$ws.Range("E9").Value = "  +3.45%  "

$ws.Range("D10").Value = "42.55"
$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("D11").Value = "1.122"
$ws.Range("E11").Value = "  +3.59%  "

$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").Value = "21.13"
$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "6.191"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.340"
$ws.Range("E15").Value = "  +5.38%  "

$ws.Range("D16").Value = "1.801.68"
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").Value = "92.35"
$ws.Range("E17").Value = "  +4.46%  "

$ws.Range("E18").Value = "  +2.48%  "

$ws.Range("D19").Value = "0.06443"
$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").Value = "0.9993"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "17.36"
$ws.Range("E21").Value = "  +4.05%  "

$ws.Range("D22").Value = "6.003"
$ws.Range("E22").Value = "  +2.67%  "

$ws.Range("D23").Value = "28.357.74"
$ws.Range("E23").Value = "  +3.70%  "

$ws.Range("D24").Value = "11.46"
$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").Value = "2.127"
$ws.Range("E25").Value = "  +3.14%  "

$ws.Range("D26").Value = "158.74"
$ws.Range("E26").Value = "  +2.83%  "

$ws.Range("D27").Value = "20.71"
$ws.Range("E27").Value = "  +3.21%  "

$ws.Range("D28").Value = "2.399"
$ws.Range("E28").Value = "  +2.75%  "

$ws.Range("D29").Value = "2.006.73"
$ws.Range("E29").Value = "  +1.81%  "

$ws.Range("D30").Value = "123.74"
$ws.Range("E30").Value = "  +2.03%  "

$ws.Range("D31").Value = "1.125"
$ws.Range("E31").Value = "  +5.98%  "

$ws.Range("E32").Value = "  +4.60%  "

$ws.Range("D33").Value = "5.747"
$ws.Range("E33").Value = "  +3.49%  "

$ws.Range("D34").Value = "3.667"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").Value = "0.2322"
$ws.Range("E35").Value = "  +14.50%  "

$ws.Range("D36").Value = "0.06451"
$ws.Range("E36").Value = "  +7.95%  "

$ws.Range("D37").Value = "0.02321"
$ws.Range("E37").Value = "  +4.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.830"
$ws.Range("E38").Value = "  +10.70%  "

$ws.Range("D39").Value = "11.65"
$ws.Range("E39").Value = "  +4.14%  "

$ws.Range("D40").Value = "5.075"
$ws.Range("E40").Value = "  +5.19%  "

$ws.Range("D41").Value = "0.6421"
$ws.Range("E41").Value = "  +4.98%  "

$ws.Range("D42").Value = "0.9993"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.160"
$ws.Range("E43").Value = "  +1.85%  "

$ws.Range("E44").Value = "  -3.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("E45").Value = "  +2.42%  "

$ws.Range("E46").Value = "  +4.01%  "

$ws.Range("E47").Value = "  +1.83%  "

$ws.Range("D48").Value = "125.24"
$ws.Range("E48").Value = "  +3.21%  "

$ws.Range("D49").Value = "1.985"
$ws.Range("E49").Value = "  +5.44%  "

$ws.Range("D51").Value = "0.06906"
$ws.Range("E51").Value = "  +3.12%  "
